$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.780.00'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '1.894.32'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.26'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4729'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2927'
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06522'
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07788'
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7418'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = '1.889.45'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.90'
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.249'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.25'
$ws.Range("E16").Value = '  +4.19%  '
$ws.Range("D17").Value = '30.762.15'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.27'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007519'
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9994'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '2.131.30'
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.329'
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9986'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.288'
$ws.Range("E24").Value = '  +2.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.233'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.47'
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.02'
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.924'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09806'
$ws.Range("E29").Value = '  -1.91%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.343'
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.490'
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.303'
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.188'
$ws.Range("E33").Value = '  +2.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04907'
$ws.Range("E34").Value = '  +2.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.135'
$ws.Range("E35").Value = '  +1.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6991'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.708'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01900'
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.842'
$ws.Range("E39").Value = '  +3.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.15'
$ws.Range("E40").Value = '  +4.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.315'
$ws.Range("E41").Value = '  +1.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.016'
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4300'
$ws.Range("E43").Value = '  +2.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8339'
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.82'
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.593'
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.038'
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.41'
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '911.64'
$ws.Range("E50").Value = '  -1.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3981'
$ws.Range("E51").Value = '  +2.94%  '
